$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date format (style) from A464 into the new date cells A465:A491
$ws.Range("A464").Copy()
$ws.Range("A465:A491").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Build the new data block (rows 465-491, columns A-D)
$data = @(
  @(44539, 0, 1, 46.70714619336758),
  @(44540, 1, 1, 46.70714619336758),
  @(44541, 0, 1, 46.70714619336758),
  @(44542, 0, 1, 46.70714619336758),
  @(44543, 0, 1, 46.70714619336758),
  @(44544, 0, 1, 46.70714619336758),
  @(44545, 0, 1, 46.70714619336758),
  @(44546, 3, 4, 186.8285847734703),
  @(44547, 0, 3, 140.1214385801028),
  @(44548, 0, 3, 140.1214385801028),
  @(44550, 0, 3, 140.1214385801028),
  @(44551, 0, 3, 140.1214385801028),
  @(44552, 0, 3, 140.1214385801028),
  @(44553, 1, 4, 186.8285847734703),
  @(44554, 0, 1, 46.70714619336758),
  @(44555, 4, 5, 233.5357309668379),
  @(44556, 4, 9, 420.3643157403083),
  @(44557, 6, 15, 700.6071929005137),
  @(44558, 3, 18, 840.7286314806165),
  @(44559, 2, 20, 934.1429238673517),
  @(44560, 5, 24, 1120.971508640822),
  @(44561, 9, 33, 1541.33582438113),
  @(44562, 9, 38, 1774.871555347968),
  @(44563, 2, 36, 1681.457262961233),
  @(44564, 3, 33, 1541.33582438113),
  @(44565, 3, 33, 1541.33582438113),
  @(44566, 2, 33, 1541.33582438113)
)

$rowCount = $data.Count
$arr = New-Object "object[,]" $rowCount,4
for ($i = 0; $i -lt $rowCount; $i++) {
  for ($j = 0; $j -lt 4; $j++) {
    $arr[$i,$j] = $data[$i][$j]
  }
}

$ws.Range("A465:D491").Value = $arr

Write-Host "Added rows 465 to 491"
